$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last row (username4), restoring the sheet to A1:A4
$ws.Range("A5").EntireRow.Delete()

# Restore the original (pre-revision) username order
$ws.Range("A1").Value = "username5"
$ws.Range("A2").Value = "username1"
$ws.Range("A3").Value = "username2"
$ws.Range("A4").Value = "username3"

# Match the saved selection state (active cell now A4, the last row)
$ws.Range("A4").Select()
